# Progress update: "update semaphore update conspect update progress"
#
# The tracker's "Анализ" (analysis) sheet pulls its percentages from small
# input tables on each topic sheet via plain cell formulas, and a bar
# chart on "Анализ" is wired to those same cells. Bumping the finished-item
# counters on the "Java" sheet is the actual authored edit here; everything
# downstream (the "Анализ" percentages/ETA math and the chart) recomputes
# from those inputs automatically on recalculation.

$wb = $excel.ActiveWorkbook

$analysis = $wb.Worksheets.Item(1)       # "Анализ" - stays the active tab
$java     = $wb.Worksheets.Item("Java")

# "Продвинутая Java (курс)" progress: 25/43 -> 27/43 completed items
$java.Range("A3").Value = 27

# "Java вопросы собеседований" progress: 12/179 -> 13/179 completed items
$java.Range("A6").Value = 13

# Remember the cursor position on the Java sheet (D16 -> D13), the way
# Excel persists each sheet's last selection, then return focus to the
# sheet that was actually active before this edit.
$java.Activate() | Out-Null
$java.Range("D13").Select() | Out-Null
$analysis.Activate() | Out-Null
